$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update changed Column D text for existing rows ---
$ws.Range("D1").Value = "Just a birthday at the squad's home. Happy Birthday, Dao."
$ws.Range("D3").Value = 'Just a nice day on December 31st 23:59PM. '
$ws.Range("D4").Value = 'If Dao is the champion versus the player, I am <a href=\"https://github.com/taylorotwell\">Taylor Otwell</a>.'
$ws.Range("D5").Value = 'The most badass wallpaper ever. I bought my Saber Hellfire because of this, setting this to be my PC wallpaper, and also ever since this was discovered, I learn fullstack webdev (useless coincidental fact for you). '
$ws.Range("D7").Value = 'Next time, use Uber, you blue little boy. You do own your own phone, am I right?'
$ws.Range("D13").Value = 'Did you know that this day honors Saint Valentine who was martyred in 269 AD?'

# --- Add new rows 16-19 ---
$ws.Range("A16").Value = 16
$ws.Range("B16").Value = '3… 2… 1… GO!'
$ws.Range("C16").Value = '16.jpg'
$ws.Range("D16").Value = 'Lodu still do the heck anywhere and its time for The Squad to end this once and for all.'

$ws.Range("A17").Value = 17
$ws.Range("B17").Value = 'Autumn Exploration'
$ws.Range("C17").Value = '17.jpg'
$ws.Range("D17").Value = 'Bazzi is sleeping, the couple is picking a bunch of apples, Uni transports them, and Ethen? Finding something interesting?'

$ws.Range("A18").Value = 18
$ws.Range("B18").Value = 'Winter for Winners'
$ws.Range("C18").Value = '18.jpg'
$ws.Range("D18").Value = 'Going outside to do a fun activities. Racing is fun, but not always be a primary source of fun. Take snowman assembling or skating for example :)'

$ws.Range("A19").Value = 19
$ws.Range("B19").Value = 'CNY Cooking Time'
$ws.Range("C19").Value = '19.jpg'
$ws.Range("D19").Value = 'Yummy… what meat is that? Cuz I am quite allergic to seafood :('

# --- Wrap text on the whole Column D data range (D1:D19) ---
$ws.Range("D1:D19").WrapText = $true

# --- Row heights: row5 is an explicit custom height; rows 4/17/18 get the
#     auto two-line wrap height from their wrapped content ---
$ws.Rows.Item(4).RowHeight = 28.8
$ws.Rows.Item(5).RowHeight = 30.6
$ws.Rows.Item(17).RowHeight = 28.8
$ws.Rows.Item(18).RowHeight = 28.8

# --- Selection / view ---
$null = $ws.Range("D18").Select()
